$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New B (id) / C (speaker_variant) values for rows 2-22.
# D (is_prefered) is cleared for every row in this export.
$rows = @(
    @{ Row = 2;  B = "#pythonissa";     C = "Pythonissa" },
    @{ Row = 3;  B = "#propheten";      C = "Propheten" },
    @{ Row = 4;  B = "#abiathar";       C = "Abiathar" },
    @{ Row = 5;  B = "#priester";       C = "Priester" },
    @{ Row = 6;  B = "#boose";          C = "Boose" },
    @{ Row = 7;  B = "#michol";         C = "Michol" },
    @{ Row = 8;  B = "#voester";        C = "Voester" },
    @{ Row = 9;  B = "#saul";           C = "Saul" },
    @{ Row = 10; B = "#abisai";         C = "Abisai" },
    @{ Row = 11; B = "#godt-den-heere"; C = "Godt den Heere" },
    @{ Row = 12; B = "#achinoam";       C = "Achinoam" },
    @{ Row = 13; B = "#ionathan";       C = "Ionathan" },
    @{ Row = 14; B = "#dauid";          C = "Dauid" },
    @{ Row = 15; B = "#echo";           C = "Echo" },
    @{ Row = 16; B = "#amalechiet";     C = "Amalechiet" },
    @{ Row = 17; B = "#samuel";         C = "Samuel" },
    @{ Row = 18; B = "#arach";          C = "Arach" },
    @{ Row = 19; B = "#abigail";        C = "Abigail" },
    @{ Row = 20; B = "#pashur";         C = "Pashur" },
    @{ Row = 21; B = "#abner";          C = "Abner" },
    @{ Row = 22; B = "#egiptenaer";     C = "Egiptenaer" }
)

foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 2).Value = $r.B
    $ws.Cells.Item($r.Row, 3).Value = $r.C
    $ws.Cells.Item($r.Row, 4).ClearContents()
}

Write-Output "applied"
